# Lesson 1 prep assignment: turn the "dice" picture + two ad-hoc numbered
# questions into a proper captioned figure.
#
#   1. The paragraph holding the dice picture switches from the
#      "BodyText" style to "CaptionedFigure".
#   2. The stray trailing " " run that followed the picture is removed.
#   3. The two numbered ("An eight sided die has..." / "Calculate the
#      probability of not...") paragraphs go away; in their place a single
#      "ImageCaption" styled paragraph reading "dice" is left behind.

$d = $word.ActiveDocument

# Find the paragraph that actually contains the inline picture, instead of
# assuming a fixed paragraph index.
$shape = $d.InlineShapes(1)
$picPara = $shape.Range.Paragraphs(1)

$picPara.Style = "CaptionedFigure"

# Drop the lone trailing space run that trails the picture inside its
# paragraph.
$shape.Range.Delete()

# The two "dice" questions are the paragraphs right after the picture.
$question1 = $picPara.Next()
$question2 = $question1.Next()

# Remove the second question entirely...
$question2.Range.Delete()

# ...and turn the first one into the new image caption.
$question1.Range.Text = "dice"
$question1.Style = "ImageCaption"
